# Append new case rows (1064-1075) for case 21TRC08418 to Sheet1.
# These rows record charges for Bunner/Hemmeter tied to case 21TRC08418,
# including the OVI counts that drop "DIP" from the jail-reporting term
# and change the jail-day warn outcome to a pass (time served = 177 of 180).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Columns: A Case#, B Name, C Charge, D Statute, E Class, F Plea,
#          G Finding, H Fine, I Fine Suspended, J Jail Days, K Jail Days Suspended
# "TEXT" cells must keep their literal/textual form (Excel would otherwise
# auto-coerce numeric-looking strings like "4511.33" or "$ 375" into numbers),
# so a "@" (Text) number format is applied before the value is written, then
# the style is reset to Normal so no lingering formatting is left behind.
# $null marks a cell that must not be created at all; "EMPTY" marks a cell
# that must exist but stay blank (matches the source workbook's
# "present-but-empty" cells).

$rows = @(
    @("21TRC08418","Bunner",  "Driving In Marked Lanes", "4511.33",     "MM","Dismissed", $null,    " ",     " ",   " ",  " "),
    @("21TRC08418","Bunner",  "Turn And Stop Signals",   "No Data",     "MM","Dismissed", $null,    " ",     " ",   " ",  " "),
    @("21TRC08418","Bunner",  "OVI Alcohol / Drugs 1st",  "4511.19A1A*","M1","No Contest","Guilty", "$ 375", "$ 0", "180","177"),
    @("21TRC08418","Hemmeter","Driving In Marked Lanes", "4511.33",     "MM","Dismissed", $null,    " ",     " ",   " ",  " "),
    @("21TRC08418","Hemmeter","Turn And Stop Signals",   "No Data",     "MM","Dismissed", $null,    " ",     " ",   " ",  " "),
    @("21TRC08418","Hemmeter","OVI Alcohol / Drugs 1st",  "4511.19A1A*","M1","No Contest","Guilty", "$ 375", "$ 0", "180","177"),
    @("21TRC08418","Hemmeter","Driving In Marked Lanes", "4511.33",     "MM","Dismissed", $null,    " ",     " ",   " ",  " "),
    @("21TRC08418","Hemmeter","Turn And Stop Signals",   "No Data",     "MM","Dismissed", $null,    " ",     " ",   " ",  " "),
    @("21TRC08418","Hemmeter","OVI Alcohol / Drugs 1st",  "4511.19A1A*","M1","No Contest","Guilty", "$ 375", "$ 0", "180","177"),
    @("21TRC08418","Hemmeter","Driving In Marked Lanes", "4511.33",     "MM","Dismissed", "EMPTY",  " ",     " ",   " ",  " "),
    @("21TRC08418","Hemmeter","Turn And Stop Signals",   "No Data",     "MM","Dismissed", "EMPTY",  " ",     " ",   " ",  " "),
    @("21TRC08418","Hemmeter","OVI Alcohol / Drugs 1st",  "4511.19A1A*","M1","No Contest","Guilty", "$ 375", "$ 0", "180","177")
)

$startRow = 1064
$textCols = @(4, 8, 9, 10, 11)   # D, H, I, J, K must stay text

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $vals = $rows[$i]

    for ($col = 1; $col -le 11; $col++) {
        $val = $vals[$col - 1]
        $cell = $ws.Cells.Item($r, $col)

        if ($col -eq 7) {
            # Finding column: $null => leave cell absent, "EMPTY" => present but blank,
            # otherwise write the literal finding text.
            if ($val -eq "EMPTY") {
                $cell.NumberFormat = "@"
                $cell.Value = ""
                $cell.Style = "Normal"
            } elseif ($val -ne $null) {
                $cell.Value = $val
            }
            continue
        }

        if ($textCols -contains $col) {
            $cell.NumberFormat = "@"
            $cell.Value = $val
            $cell.Style = "Normal"
        } else {
            $cell.Value = $val
        }
    }
}
